$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cell, $text) {
    # Force the numeric-looking string to be stored as text (matching the
    # original inline-string cell type), then restore the default "Normal"
    # style so no stray style index is left on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '58.788.70'
Set-TextValue $ws.Cells.Item(2, 5) '  -3.22%  '
Set-TextValue $ws.Cells.Item(3, 4) '2.554.06'
Set-TextValue $ws.Cells.Item(3, 5) '  -1.67%  '
Set-TextValue $ws.Cells.Item(4, 4) '1.00'
Set-TextValue $ws.Cells.Item(4, 5) '  -0.07%  '
Set-TextValue $ws.Cells.Item(5, 4) '504.69'
Set-TextValue $ws.Cells.Item(5, 5) '  -3.55%  '
Set-TextValue $ws.Cells.Item(6, 4) '141.79'
Set-TextValue $ws.Cells.Item(6, 5) '  -8.00%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.999'
Set-TextValue $ws.Cells.Item(7, 5) '  +0.00%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.551'
Set-TextValue $ws.Cells.Item(8, 5) '  -6.04%  '
Set-TextValue $ws.Cells.Item(9, 4) '2.555.20'
Set-TextValue $ws.Cells.Item(9, 5) '  -1.79%  '
Set-TextValue $ws.Cells.Item(10, 4) '6.20'
Set-TextValue $ws.Cells.Item(10, 5) '  -7.08%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.101'
Set-TextValue $ws.Cells.Item(11, 5) '  -4.32%  '
Set-TextValue $ws.Cells.Item(12, 4) '0.328'
Set-TextValue $ws.Cells.Item(12, 5) '  -4.96%  '
Set-TextValue $ws.Cells.Item(13, 4) '0.128'
Set-TextValue $ws.Cells.Item(13, 5) '  -0.96%  '
Set-TextValue $ws.Cells.Item(14, 4) '3.001.38'
Set-TextValue $ws.Cells.Item(14, 5) '  -1.65%  '
Set-TextValue $ws.Cells.Item(15, 4) '58.741.53'
Set-TextValue $ws.Cells.Item(15, 5) '  -3.34%  '
Set-TextValue $ws.Cells.Item(16, 4) '20.45'
Set-TextValue $ws.Cells.Item(16, 5) '  -5.09%  '
Set-TextValue $ws.Cells.Item(17, 4) '0.0000134'
Set-TextValue $ws.Cells.Item(17, 5) '  -5.02%  '
Set-TextValue $ws.Cells.Item(18, 4) '2.561.58'
Set-TextValue $ws.Cells.Item(18, 5) '  -1.55%  '
Set-TextValue $ws.Cells.Item(19, 4) '4.49'
Set-TextValue $ws.Cells.Item(19, 5) '  -5.62%  '
Set-TextValue $ws.Cells.Item(20, 4) '329.77'
Set-TextValue $ws.Cells.Item(20, 5) '  -7.07%  '
Set-TextValue $ws.Cells.Item(21, 4) '10.00'
Set-TextValue $ws.Cells.Item(21, 5) '  -5.11%  '
Set-TextValue $ws.Cells.Item(22, 4) '1.00'
Set-TextValue $ws.Cells.Item(22, 5) '  +0.23%  '
Set-TextValue $ws.Cells.Item(23, 4) '5.90'
Set-TextValue $ws.Cells.Item(23, 5) '  -4.55%  '
Set-TextValue $ws.Cells.Item(24, 4) '59.32'
Set-TextValue $ws.Cells.Item(24, 5) '  -2.86%  '
Set-TextValue $ws.Cells.Item(25, 4) '0.403'
Set-TextValue $ws.Cells.Item(25, 5) '  -5.14%  '
Set-TextValue $ws.Cells.Item(26, 4) '1.00'
Set-TextValue $ws.Cells.Item(26, 5) '  -0.02%  '
Set-TextValue $ws.Cells.Item(27, 4) '0.157'
Set-TextValue $ws.Cells.Item(27, 5) '  -5.65%  '
Set-TextValue $ws.Cells.Item(28, 4) '0.0₃0771'
Set-TextValue $ws.Cells.Item(28, 5) '  -8.23%  '
Set-TextValue $ws.Cells.Item(29, 4) '6.83'
Set-TextValue $ws.Cells.Item(29, 5) '  -7.36%  '
Set-TextValue $ws.Cells.Item(30, 4) '1.00'
Set-TextValue $ws.Cells.Item(30, 5) '  +0.02%  '
Set-TextValue $ws.Cells.Item(31, 4) '149.30'
Set-TextValue $ws.Cells.Item(31, 5) '  -0.26%  '
Set-TextValue $ws.Cells.Item(32, 4) '18.45'
Set-TextValue $ws.Cells.Item(32, 5) '  -4.80%  '
Set-TextValue $ws.Cells.Item(33, 4) '1.53'
Set-TextValue $ws.Cells.Item(33, 5) '  -4.16%  '
Set-TextValue $ws.Cells.Item(34, 4) '5.76'
Set-TextValue $ws.Cells.Item(34, 5) '  -8.16%  '
Set-TextValue $ws.Cells.Item(35, 4) '3.87'
Set-TextValue $ws.Cells.Item(35, 5) '  -7.55%  '
Set-TextValue $ws.Cells.Item(36, 4) '0.866'
Set-TextValue $ws.Cells.Item(36, 5) '  -5.46%  '
Set-TextValue $ws.Cells.Item(37, 4) '1.09'
Set-TextValue $ws.Cells.Item(37, 5) '  -8.28%  '
Set-TextValue $ws.Cells.Item(38, 4) '35.76'
Set-TextValue $ws.Cells.Item(38, 5) '  -1.77%  '
Set-TextValue $ws.Cells.Item(39, 4) '0.819'
Set-TextValue $ws.Cells.Item(39, 5) '  -9.87%  '
Set-TextValue $ws.Cells.Item(40, 4) '284.86'
Set-TextValue $ws.Cells.Item(40, 5) '  -2.25%  '
Set-TextValue $ws.Cells.Item(41, 4) '1.37'
Set-TextValue $ws.Cells.Item(41, 5) '  -7.87%  '
Set-TextValue $ws.Cells.Item(42, 4) '3.48'
Set-TextValue $ws.Cells.Item(42, 5) '  -8.12%  '
Set-TextValue $ws.Cells.Item(43, 4) '0.998'
Set-TextValue $ws.Cells.Item(43, 5) '  -0.01%  '
Set-TextValue $ws.Cells.Item(44, 4) '0.0978'
Set-TextValue $ws.Cells.Item(44, 5) '  -3.33%  '
Set-TextValue $ws.Cells.Item(45, 4) '0.603'
Set-TextValue $ws.Cells.Item(45, 5) '  -3.05%  '
Set-TextValue $ws.Cells.Item(46, 4) '0.0527'
Set-TextValue $ws.Cells.Item(46, 5) '  -5.44%  '
Set-TextValue $ws.Cells.Item(47, 4) '10.34'
Set-TextValue $ws.Cells.Item(47, 5) '  +0.04%  '
Set-TextValue $ws.Cells.Item(48, 4) '18.52'
Set-TextValue $ws.Cells.Item(48, 5) '  -5.14%  '
Set-TextValue $ws.Cells.Item(49, 4) '0.0225'
Set-TextValue $ws.Cells.Item(49, 5) '  -5.31%  '
Set-TextValue $ws.Cells.Item(50, 4) '4.50'
Set-TextValue $ws.Cells.Item(50, 5) '  -8.19%  '
Set-TextValue $ws.Cells.Item(51, 4) '1.878.67'
Set-TextValue $ws.Cells.Item(51, 5) '  -4.27%  '
